# This workbook's data rows (2-8) get cyclically rotated:
#   new row 2 <- old row 4
#   new row 3 <- old row 5
#   new row 4 <- old row 6
#   new row 5 <- old row 7
#   new row 6 <- old row 8
#   new row 7 <- old row 2
#   new row 8 <- old row 3
#
# We stage the rotated rows in a scratch area far below the used range
# (rows 202-208). The set of populated columns for each row is determined
# once, up front, directly from the untouched source rows (2-8) - this is
# important because Range.Copy of a genuinely-empty cell produces a
# "typeless" destination cell whose Value2 reads back as $null, which
# would make a second hole-detection pass (run against already-copied
# scratch rows) unreliable. By recording the column layout up front we
# avoid ever needing to re-detect holes in data that has already been
# copied once.
#
# Range.Copy (rather than assigning .Value/.Value2 with literal strings)
# is used throughout so that cell contents are carried over exactly as
# stored, without Excel re-interpreting date-like text as real dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstCol = 1   # A
$lastCol = 51   # AY

function ColLetters($n) {
    $letters = ""
    while ($n -gt 0) {
        $rem = ($n - 1) % 26
        $letters = [char](65 + $rem) + $letters
        $n = [int](($n - $rem - 1) / 26)
    }
    return $letters
}

# Returns a list of "StartCol:EndCol" (letters) contiguous runs of
# populated cells on the given row.
function GetRuns($row) {
    $runs = New-Object System.Collections.ArrayList
    $col = $firstCol
    while ($col -le $lastCol) {
        $letters = ColLetters($col)
        $v = $ws.Range($letters + $row).Value2
        if ($v -eq $null) {
            $col = $col + 1
            continue
        }
        $runStart = $col
        $runEnd = $col
        $nextCol = $col + 1
        while ($nextCol -le $lastCol) {
            $nextRef = (ColLetters($nextCol)) + $row
            $nv = $ws.Range($nextRef).Value2
            if ($nv -eq $null) {
                break
            }
            $runEnd = $nextCol
            $nextCol = $nextCol + 1
        }
        [void]$runs.Add(@{ Start = (ColLetters($runStart)); End = (ColLetters($runEnd)) })
        $col = $runEnd + 1
    }
    return $runs
}

function CopyRowUsingRuns($srcRow, $dstRow, $runs) {
    foreach ($run in $runs) {
        $srcRange = $ws.Range($run.Start + $srcRow + ":" + $run.End + $srcRow)
        $dstRange = $ws.Range($run.Start + $dstRow + ":" + $run.End + $dstRow)
        $srcRange.Copy($dstRange)
    }
}

$mapping = @{ 2 = 4; 3 = 5; 4 = 6; 5 = 7; 6 = 8; 7 = 2; 8 = 3 }
$scratchOffset = 200

# Determine the populated-column layout of every source row up front,
# while rows 2-8 are still in their pristine, unmodified state.
$rowRuns = @{}
foreach ($row in 2..8) {
    $rowRuns[$row] = GetRuns($row)
}

# Step 1: copy each source row into a scratch row (scratchOffset + destination row)
foreach ($destRow in 2..8) {
    $srcRow = $mapping[$destRow]
    $scratchRow = $scratchOffset + $destRow
    CopyRowUsingRuns $srcRow $scratchRow $rowRuns[$srcRow]
}

# Step 2: clear the original rows 2-8 entirely
$ws.Range("A2:AY8").ClearContents()

# Step 3: copy the scratch rows back into rows 2-8, using the same known
# column layout as the original source row (content and layout travel
# together).
foreach ($destRow in 2..8) {
    $srcRow = $mapping[$destRow]
    $scratchRow = $scratchOffset + $destRow
    CopyRowUsingRuns $scratchRow $destRow $rowRuns[$srcRow]
}

# Step 4: clear the scratch area
$ws.Range("A202:AY208").ClearContents()
